$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 58
$ws.Cells.Item(58, 8).Value = 2019.6364  # H58: 2034.75 -> 2019.6364
$ws.Cells.Item(58, 9).Value = 29  # I58: 48.75 -> 29
$ws.Cells.Item(58, 10).Value = 3157.1428  # J58: 2531.25 -> 3157.1428
$ws.Cells.Item(58, 11).Value = 87  # K58: 146.25 -> 87
$ws.Cells.Item(58, 12).Value = 9471.428400000001  # L58: 7593.75 -> 9471.428400000001
$ws.Cells.Item(58, 13).Value = 63  # M58: 3.75 -> 63
$ws.Cells.Item(58, 14).Value = -9771.428400000001  # N58: -7893.75 -> -9771.428400000001

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Cells.Item(32, 8).Value = 28873.723  # H32: 28867.889 -> 28873.723
$ws.Cells.Item(32, 9).Value = 5184.023  # I32: 5084.0444 -> 5184.023
$ws.Cells.Item(32, 10).Value = 133108.4  # J32: 147787.11 -> 133108.4
$ws.Cells.Item(32, 11).Value = 5184.023  # K32: 5084.0444 -> 5184.023
$ws.Cells.Item(32, 12).Value = 133108.4  # L32: 147787.11 -> 133108.4
$ws.Cells.Item(32, 13).Value = -4897.023  # M32: -4797.0444 -> -4897.023
$ws.Cells.Item(32, 14).Value = -133682.4  # N32: -148361.11 -> -133682.4

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31, 8).Value = 20835438  # H31: 18869882 -> 20835438
$ws.Cells.Item(31, 9).Value = 58825180  # I31: 47620496 -> 58825180
$ws.Cells.Item(31, 10).Value = 2351.8386  # J31: 2290.8125 -> 2351.8386
$ws.Cells.Item(31, 11).Value = 58825180  # K31: 47620496 -> 58825180
$ws.Cells.Item(31, 12).Value = 2351.8386  # L31: 2290.8125 -> 2351.8386
$ws.Cells.Item(31, 13).Value = -58824885  # M31: -47620201 -> -58824885
$ws.Cells.Item(31, 14).Value = -2941.8386  # N31: -2880.8125 -> -2941.8386
# Row 34
$ws.Cells.Item(34, 8).Value = 20835438  # H34: 18869882 -> 20835438
$ws.Cells.Item(34, 9).Value = 58825180  # I34: 47620496 -> 58825180
$ws.Cells.Item(34, 10).Value = 2351.8386  # J34: 2290.8125 -> 2351.8386
$ws.Cells.Item(34, 11).Value = 58825180  # K34: 47620496 -> 58825180
$ws.Cells.Item(34, 12).Value = 2351.8386  # L34: 2290.8125 -> 2351.8386
$ws.Cells.Item(34, 13).Value = -58824978  # M34: -47620294 -> -58824978
$ws.Cells.Item(34, 14).Value = -2755.8386  # N34: -2694.8125 -> -2755.8386
# Row 99
$ws.Cells.Item(99, 8).Value = 1419.5714  # H99: 1493.2307 -> 1419.5714
$ws.Cells.Item(99, 9).Value = 1306.7273  # I99: 1384.6666 -> 1306.7273
$ws.Cells.Item(99, 10).Value = 1833.3334  # J99: 1737.5 -> 1833.3334
$ws.Cells.Item(99, 11).Value = 1306.7273  # K99: 1384.6666 -> 1306.7273
$ws.Cells.Item(99, 12).Value = 1833.3334  # L99: 1737.5 -> 1833.3334
$ws.Cells.Item(99, 13).Value = 191.2727  # M99: 113.3334 -> 191.2727
$ws.Cells.Item(99, 14).Value = -4829.3334  # N99: -4733.5 -> -4829.3334
# Row 109
$ws.Cells.Item(109, 8).Value = 21156.6  # H109: 21286.6 -> 21156.6
$ws.Cells.Item(109, 10).Value = 21156.6  # J109: 21286.6 -> 21156.6
$ws.Cells.Item(109, 12).Value = 21156.6  # L109: 21286.6 -> 21156.6
$ws.Cells.Item(109, 14).Value = -23236.6  # N109: -23366.6 -> -23236.6
# Row 126
$ws.Cells.Item(126, 8).Value = 1419.5714  # H126: 1493.2307 -> 1419.5714
$ws.Cells.Item(126, 9).Value = 1306.7273  # I126: 1384.6666 -> 1306.7273
$ws.Cells.Item(126, 10).Value = 1833.3334  # J126: 1737.5 -> 1833.3334
$ws.Cells.Item(126, 11).Value = 3920.1819  # K126: 4153.9998 -> 3920.1819
$ws.Cells.Item(126, 12).Value = 5500.0002  # L126: 5212.5 -> 5500.0002
$ws.Cells.Item(126, 13).Value = -1450.1819  # M126: -1683.9998 -> -1450.1819
$ws.Cells.Item(126, 14).Value = -10440.0002  # N126: -10152.5 -> -10440.0002

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 11
$ws.Cells.Item(11, 8).Value = 849.2941  # H11: 703.625 -> 849.2941
$ws.Cells.Item(11, 9).Value = 39.833332  # I11: 66 -> 39.833332
$ws.Cells.Item(11, 10).Value = 1290.8182  # J11: 1086.2 -> 1290.8182
$ws.Cells.Item(11, 11).Value = 119.499996  # K11: 198 -> 119.499996
$ws.Cells.Item(11, 12).Value = 3872.4546  # L11: 3258.6 -> 3872.4546
$ws.Cells.Item(11, 13).Value = 20.500004  # M11: -58 -> 20.500004
$ws.Cells.Item(11, 14).Value = -4152.4546  # N11: -3538.6 -> -4152.4546
# Row 68
$ws.Cells.Item(68, 8).Value = 1311.0613  # H68: 1344.0444 -> 1311.0613
$ws.Cells.Item(68, 9).Value = 969.05884  # I68: 944.93335 -> 969.05884
$ws.Cells.Item(68, 10).Value = 1492.75  # J68: 1543.6 -> 1492.75
$ws.Cells.Item(68, 11).Value = 2907.17652  # K68: 2834.80005 -> 2907.17652
$ws.Cells.Item(68, 12).Value = 4478.25  # L68: 4630.799999999999 -> 4478.25
$ws.Cells.Item(68, 13).Value = -2096.17652  # M68: -2023.80005 -> -2096.17652
$ws.Cells.Item(68, 14).Value = -6100.25  # N68: -6252.799999999999 -> -6100.25
# Row 70
$ws.Cells.Item(70, 8).Value = 5568.636  # H70: 5143.875 -> 5568.636
$ws.Cells.Item(70, 9).Value = 1952.75  # I70: 2059.1428 -> 1952.75
$ws.Cells.Item(70, 10).Value = 7634.857  # J70: 7543.1113 -> 7634.857
$ws.Cells.Item(70, 11).Value = 5858.25  # K70: 6177.428400000001 -> 5858.25
$ws.Cells.Item(70, 12).Value = 22904.571  # L70: 22629.3339 -> 22904.571
$ws.Cells.Item(70, 13).Value = -5543.25  # M70: -5862.428400000001 -> -5543.25
$ws.Cells.Item(70, 14).Value = -23534.571  # N70: -23259.3339 -> -23534.571
# Row 71
$ws.Cells.Item(71, 8).Value = 1311.0613  # H71: 1344.0444 -> 1311.0613
$ws.Cells.Item(71, 9).Value = 969.05884  # I71: 944.93335 -> 969.05884
$ws.Cells.Item(71, 10).Value = 1492.75  # J71: 1543.6 -> 1492.75
$ws.Cells.Item(71, 11).Value = 8721.529560000001  # K71: 8504.400149999999 -> 8721.529560000001
$ws.Cells.Item(71, 12).Value = 13434.75  # L71: 13892.4 -> 13434.75
$ws.Cells.Item(71, 13).Value = -4665.529560000001  # M71: -4448.400149999999 -> -4665.529560000001
$ws.Cells.Item(71, 14).Value = -21546.75  # N71: -22004.4 -> -21546.75
# Row 73
$ws.Cells.Item(73, 8).Value = 5568.636  # H73: 5143.875 -> 5568.636
$ws.Cells.Item(73, 9).Value = 1952.75  # I73: 2059.1428 -> 1952.75
$ws.Cells.Item(73, 10).Value = 7634.857  # J73: 7543.1113 -> 7634.857
$ws.Cells.Item(73, 11).Value = 5858.25  # K73: 6177.428400000001 -> 5858.25
$ws.Cells.Item(73, 12).Value = 22904.571  # L73: 22629.3339 -> 22904.571
$ws.Cells.Item(73, 13).Value = -4766.25  # M73: -5085.428400000001 -> -4766.25
$ws.Cells.Item(73, 14).Value = -25088.571  # N73: -24813.3339 -> -25088.571
# Row 74
$ws.Cells.Item(74, 8).Value = 6714.2856  # H74: 8311.429 -> 6714.2856
$ws.Cells.Item(74, 10).Value = 9000  # J74: 9530 -> 9000
$ws.Cells.Item(74, 12).Value = 27000  # L74: 28590 -> 27000
$ws.Cells.Item(74, 14).Value = -29122  # N74: -30712 -> -29122
# Row 77
$ws.Cells.Item(77, 8).Value = 6714.2856  # H77: 8311.429 -> 6714.2856
$ws.Cells.Item(77, 10).Value = 9000  # J77: 9530 -> 9000
$ws.Cells.Item(77, 12).Value = 81000  # L77: 85770 -> 81000
$ws.Cells.Item(77, 14).Value = -91608  # N77: -96378 -> -91608
# Row 131
$ws.Cells.Item(131, 8).Value = 898.0303  # H131: 923.05554 -> 898.0303
$ws.Cells.Item(131, 9).Value = 357.14285  # I131: 366.66666 -> 357.14285
$ws.Cells.Item(131, 10).Value = 1043.6538  # J131: 1034.3334 -> 1043.6538
$ws.Cells.Item(131, 11).Value = 1071.42855  # K131: 1099.99998 -> 1071.42855
$ws.Cells.Item(131, 12).Value = 3130.9614  # L131: 3103.0002 -> 3130.9614
$ws.Cells.Item(131, 13).Value = 3968.57145  # M131: 3940.00002 -> 3968.57145
$ws.Cells.Item(131, 14).Value = -13210.9614  # N131: -13183.0002 -> -13210.9614

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 21
$ws.Cells.Item(21, 8).Value = 750  # H21: 1049.75 -> 750
$ws.Cells.Item(21, 9).Value = 750  # I21: 900 -> 750
$ws.Cells.Item(21, 10).Value = 0  # J21: 1199.5 -> 0
$ws.Cells.Item(21, 11).Value = 750  # K21: 900 -> 750
$ws.Cells.Item(21, 12).Value = 0  # L21: 1199.5 -> 0
$ws.Cells.Item(21, 13).Value = -577  # M21: -727 -> -577
$ws.Cells.Item(21, 14).ClearContents()  # N21 removed (was -1545.5)
# Row 30
$ws.Cells.Item(30, 8).Value = 750  # H30: 1049.75 -> 750
$ws.Cells.Item(30, 9).Value = 750  # I30: 900 -> 750
$ws.Cells.Item(30, 10).Value = 0  # J30: 1199.5 -> 0
$ws.Cells.Item(30, 11).Value = 750  # K30: 900 -> 750
$ws.Cells.Item(30, 12).Value = 0  # L30: 1199.5 -> 0
$ws.Cells.Item(30, 13).Value = -645  # M30: -795 -> -645
$ws.Cells.Item(30, 14).ClearContents()  # N30 removed (was -1409.5)
# Row 57
$ws.Cells.Item(57, 8).Value = 21335.166  # H57: 21459.385 -> 21335.166
$ws.Cells.Item(57, 10).Value = 23002  # J57: 22997.666 -> 23002
$ws.Cells.Item(57, 12).Value = 23002  # L57: 22997.666 -> 23002
$ws.Cells.Item(57, 14).Value = -24642  # N57: -24637.666 -> -24642
# Row 69
$ws.Cells.Item(69, 8).Value = 150000  # H69: 0 -> 150000
$ws.Cells.Item(69, 10).Value = 150000  # J69: 0 -> 150000
$ws.Cells.Item(69, 12).Value = 150000  # L69: 0 -> 150000
$ws.Cells.Item(69, 14).Value = -151498  # N69 new cell
# Row 72
$ws.Cells.Item(72, 8).Value = 150000  # H72: 0 -> 150000
$ws.Cells.Item(72, 10).Value = 150000  # J72: 0 -> 150000
$ws.Cells.Item(72, 12).Value = 450000  # L72: 0 -> 450000
$ws.Cells.Item(72, 14).Value = -457488  # N72 new cell
# Row 126
$ws.Cells.Item(126, 8).Value = 4626.207  # H126: 5820 -> 4626.207
$ws.Cells.Item(126, 9).Value = 2498  # I126: 3000 -> 2498
$ws.Cells.Item(126, 10).Value = 9355.556  # J126: 8640 -> 9355.556
$ws.Cells.Item(126, 11).Value = 7494  # K126: 9000 -> 7494
$ws.Cells.Item(126, 12).Value = 28066.668  # L126: 25920 -> 28066.668
$ws.Cells.Item(126, 13).Value = -5024  # M126: -6530 -> -5024
$ws.Cells.Item(126, 14).Value = -33006.66800000001  # N126: -30860 -> -33006.66800000001

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Cells.Item(7, 8).Value = 6668613.5  # H7: 4349595.5 -> 6668613.5
$ws.Cells.Item(7, 9).Value = 16668501  # I7: 6668227 -> 16668501
$ws.Cells.Item(7, 10).Value = 2021.7778  # J7: 2162 -> 2021.7778
$ws.Cells.Item(7, 11).Value = 16668501  # K7: 6668227 -> 16668501
$ws.Cells.Item(7, 12).Value = 2021.7778  # L7: 2162 -> 2021.7778
$ws.Cells.Item(7, 13).Value = -16668389  # M7: -6668115 -> -16668389
$ws.Cells.Item(7, 14).Value = -2245.7778  # N7: -2386 -> -2245.7778
# Row 81
$ws.Cells.Item(81, 8).Value = 28485.334  # H81: 31000 -> 28485.334
$ws.Cells.Item(81, 10).Value = 28485.334  # J81: 31000 -> 28485.334
$ws.Cells.Item(81, 12).Value = 28485.334  # L81: 31000 -> 28485.334
$ws.Cells.Item(81, 14).Value = -30481.334  # N81: -32996 -> -30481.334
# Row 84
$ws.Cells.Item(84, 8).Value = 28485.334  # H84: 31000 -> 28485.334
$ws.Cells.Item(84, 10).Value = 28485.334  # J84: 31000 -> 28485.334
$ws.Cells.Item(84, 12).Value = 85456.00199999999  # L84: 93000 -> 85456.00199999999
$ws.Cells.Item(84, 14).Value = -95440.00199999999  # N84: -102984 -> -95440.00199999999
# Row 126
$ws.Cells.Item(126, 8).Value = 6668613.5  # H126: 4349595.5 -> 6668613.5
$ws.Cells.Item(126, 9).Value = 16668501  # I126: 6668227 -> 16668501
$ws.Cells.Item(126, 10).Value = 2021.7778  # J126: 2162 -> 2021.7778
$ws.Cells.Item(126, 11).Value = 50005503  # K126: 20004681 -> 50005503
$ws.Cells.Item(126, 12).Value = 6065.3334  # L126: 6486 -> 6065.3334
$ws.Cells.Item(126, 13).Value = -50003033  # M126: -20002211 -> -50003033
$ws.Cells.Item(126, 14).Value = -11005.3334  # N126: -11426 -> -11005.3334
# Row 132
$ws.Cells.Item(132, 8).Value = 3899.6  # H132: 3667.1667 -> 3899.6
$ws.Cells.Item(132, 10).Value = 5999  # J132: 4250.75 -> 5999
$ws.Cells.Item(132, 12).Value = 17997  # L132: 12752.25 -> 17997
$ws.Cells.Item(132, 14).Value = -23057  # N132: -17812.25 -> -23057

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 109
$ws.Cells.Item(109, 8).Value = 30970.8  # H109: 31030.8 -> 30970.8
$ws.Cells.Item(109, 10).Value = 30970.8  # J109: 31030.8 -> 30970.8
$ws.Cells.Item(109, 12).Value = 30970.8  # L109: 31030.8 -> 30970.8
$ws.Cells.Item(109, 14).Value = -33744.8  # N109: -33804.8 -> -33744.8
